$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.201.49'
$ws.Range('E2').Value = '  -1.18%  '
$ws.Range('D3').Value = '1.894.90'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '245.33'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').Value = '0.684'
$ws.Range('E6').Value = '  +8.02%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '40.37'
$ws.Range('E8').Value = '  -4.89%  '
$ws.Range('D9').Value = '0.345'
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('D10').Value = '53.10'
$ws.Range('E10').Value = '  +11.19%  '
$ws.Range('D11').Value = '0.0717'
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('D12').Value = '0.0989'
$ws.Range('E12').Value = '  -0.79%  '
$ws.Range('D13').Value = '2.169.82'
$ws.Range('E13').Value = '  -0.80%  '
$ws.Range('D14').Value = '12.47'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('E15').Value = '  +0.77%  '
$ws.Range('D16').Value = '1.893.13'
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('D17').Value = '4.77'
$ws.Range('E17').Value = '  -1.80%  '
$ws.Range('D18').Value = '35.201.17'
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('D19').Value = '71.89'
$ws.Range('E19').Value = '  -0.32%  '
$ws.Range('D20').Value = '0.0₃0813'
$ws.Range('E20').Value = '  +0.26%  '
$ws.Range('D21').Value = '239.93'
$ws.Range('E21').Value = '  -1.72%  '
$ws.Range('D22').Value = '12.48'
$ws.Range('E22').Value = '  -0.45%  '
$ws.Range('D23').Value = '4.75'
$ws.Range('E23').Value = '  -3.70%  '
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('E25').Value = '  +1.33%  '
$ws.Range('D26').Value = '2.31'
$ws.Range('E26').Value = '  +9.12%  '
$ws.Range('D27').Value = '167.75'
$ws.Range('E27').Value = '  -2.16%  '
$ws.Range('D28').Value = '8.50'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '18.15'
$ws.Range('E29').Value = '  +0.71%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = '0.129'
$ws.Range('E30').Value = '  +2.78%  '
$ws.Range('D32').Value = '4.13'
$ws.Range('E32').Value = '  +0.60%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.0562'
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('B34').Value = 'BinanceUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D34').Value = '1.01'
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('E35').Value = '  +5.19%  '
$ws.Range('D36').Value = '4.07'
$ws.Range('E36').Value = '  -1.79%  '
$ws.Range('D37').Value = '0.901'
$ws.Range('E37').Value = '  -4.82%  '
$ws.Range('D38').Value = '1.48'
$ws.Range('E38').Value = '  +12.07%  '
$ws.Range('D39').Value = '2.00'
$ws.Range('E39').Value = '  -2.27%  '
$ws.Range('E40').Value = '  +0.71%  '
$ws.Range('E41').Value = '  -2.71%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').Value = '15.96'
$ws.Range('E42').Value = '  +4.71%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = '0.0629'
$ws.Range('E43').Value = '  +5.91%  '
$ws.Range('D44').Value = '89.10'
$ws.Range('E44').Value = '  -3.00%  '
$ws.Range('D45').Value = '1.339.90'
$ws.Range('E45').Value = '  -1.71%  '
$ws.Range('D46').Value = '2.39'
$ws.Range('E46').Value = '  +1.52%  '
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('D49').Value = '44.94'
$ws.Range('E49').Value = '  -7.63%  '
$ws.Range('B50').Value = 'Gas'
$ws.Range('C50').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D50').Value = '12.16'
$ws.Range('E50').Value = '  -7.15%  '
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D51').Value = '6.42'
$ws.Range('E51').Value = '  -4.04%  '
